$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.491.63"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "3.855.53"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.60"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.59%  "
$ws.Range("D7").Value = "3.855.63"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +0.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.168"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.55%  "
$ws.Range("E12").Value = "  +1.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000284"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +14.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "4.504.29"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "3.845.40"
$ws.Range("E16").Value = "  -0.61%  "
$ws.Range("D17").Value = "68.535.42"
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.49%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.112"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.70%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "471.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.52%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.735"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000159"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.98%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.64%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.23"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.61"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.16%  "
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("D31").Value = "4.009.13"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "31.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").Value = "3.821.51"
$ws.Range("E36").Value = "  -0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.92"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +20.71%  "
$ws.Range("E38").Value = "  +1.52%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.02"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.59%  "
$ws.Range("B40").Value = "Filecoin"
$ws.Range("C40").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.99"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("E41").Value = "  +0.14%  "
$ws.Range("E42").Value = "  +0.02%  "
$ws.Range("E43").Value = "  +2.78%  "
$ws.Range("E44").Value = "  +0.91%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.77"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "420.94"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.75%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000294"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "46.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.18"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0360"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.81%  "
